$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Cell -> new text value, scraped from the "Updated cryptos list" diff.
# All values must remain *text*, matching how the sheet already stores them
# (Price cells use dotted thousand separators like "68.143.24" and the
# Volume(1h) cells are space-padded percents like "  -1.39%  "). Plain
# decimal-looking prices (e.g. "0.999", "1.00") get auto-coerced to numbers
# by Excel on plain assignment, so those are written through a temporary
# Text ("@") number format and then ClearFormats() to drop that format again
# and land back on the default (unstyled) cell, same as every other value.
$updates = [ordered]@{
    'D2' = '68.143.24'
    'E2' = '  -1.39%  '
    'D3' = '3.877.42'
    'E3' = '  -1.59%  '
    'D4' = '0.999'
    'E4' = '  -0.13%  '
    'D5' = '600.09'
    'E5' = '  -0.83%  '
    'D6' = '171.52'
    'E6' = '  +1.95%  '
    'D7' = '3.878.05'
    'E7' = '  -1.68%  '
    'E8' = '  -0.01%  '
    'D9' = '0.530'
    'E9' = '  -1.04%  '
    'D10' = '0.164'
    'E10' = '  -5.24%  '
    'D11' = '6.41'
    'E11' = '  -1.20%  '
    'D12' = '0.457'
    'E12' = '  -1.90%  '
    'D13' = '0.0000259'
    'E13' = '  +0.23%  '
    'D14' = '37.07'
    'E14' = '  -1.70%  '
    'D15' = '4.525.71'
    'E15' = '  -1.60%  '
    'D16' = '3.878.52'
    'E16' = '  -1.95%  '
    'D17' = '68.235.71'
    'E17' = '  -1.24%  '
    'D18' = '18.16'
    'E18' = '  +4.21%  '
    'D19' = '7.35'
    'E19' = '  -2.17%  '
    'E20' = '  -0.27%  '
    'D21' = '10.76'
    'E21' = '  -2.07%  '
    'D22' = '466.63'
    'E22' = '  -5.97%  '
    'D23' = '0.741'
    'E23' = '  +1.13%  '
    'E24' = '  -5.79%  '
    'D25' = '83.19'
    'E25' = '  -2.35%  '
    'D26' = '2.24'
    'E26' = '  -1.96%  '
    'D27' = '12.06'
    'E27' = '  -1.20%  '
    'D28' = '1.00'
    'E28' = '  -0.05%  '
    'D29' = '9.99'
    'E29' = '  -2.82%  '
    'D30' = '2.96'
    'E30' = '  -1.22%  '
    'D31' = '4.024.28'
    'E31' = '  -1.69%  '
    'D32' = '7.78'
    'E32' = '  -0.31%  '
    'D33' = '2.32'
    'E33' = '  -3.24%  '
    'D34' = '31.22'
    'E34' = '  -2.76%  '
    'D35' = '9.51'
    'E35' = '  -0.38%  '
    'D36' = '3.835.95'
    'E36' = '  -1.81%  '
    'D37' = '0.105'
    'E37' = '  -2.85%  '
    'D38' = '3.75'
    'E38' = '  +12.94%  '
    'D39' = '1.03'
    'E39' = '  -1.87%  '
    'E40' = '  +0.76%  '
    'D41' = '5.91'
    'E41' = '  -1.73%  '
    'D42' = '1.00'
    'E42' = '  +0.00%  '
    'D43' = '0.313'
    'E43' = '  -2.94%  '
    'D44' = '0.000298'
    'E44' = '  +4.14%  '
    'D45' = '424.37'
    'E45' = '  -2.69%  '
    'D46' = '1.98'
    'E46' = '  -1.73%  '
    'E47' = '  -0.02%  '
    'B48' = 'Cosmos'
    'C48' = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
    'D48' = '8.63'
    'E48' = '  +0.02%  '
    'B49' = 'OKB'
    'C49' = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
    'D49' = '47.18'
    'E49' = '  -1.82%  '
    'D50' = '26.86'
    'E50' = '  +4.04%  '
    'D51' = '143.13'
    'E51' = '  -0.05%  '
}

foreach ($addr in $updates.Keys) {
    $value = $updates[$addr]
    $range = $ws.Range($addr)
    if ($value.Trim() -match '^[+-]?[0-9]+(\.[0-9]+)?([eE][+-]?[0-9]+)?$') {
        # Looks like a plain number -> force text so it is not coerced.
        $range.NumberFormat = "@"
        $range.Value = $value
        $range.ClearFormats()
    } else {
        $range.Value = $value
    }
}
